$d = $word.ActiveDocument

# --- Remove duplicated / stray section-heading paragraphs -----------------
# The notebook-update code used to write each day's section headings
# ("Goals", "Code Executions", "To-Do List") twice, and also emitted a
# stray "Accomplished" heading that the rest of that day's entry never
# used. Moving the notebook-update logic into its own function fixed the
# double-invocation bug. Walk the paragraphs back-to-front (so deleting
# one never shifts the index of a paragraph still to be inspected) and
# drop the extra occurrences.

$count = $d.Paragraphs.Count
$seenGoals = 0
$seenCodeExecutions = 0
$seenToDoList = 0

for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Style.NameLocal -ne "Heading 1") { continue }

    $t = $p.Range.Text.Trim()

    if ($t -eq "Goals") {
        $seenGoals = $seenGoals + 1
        # Scanning back-to-front, the first hit is the trailing duplicate.
        if ($seenGoals -eq 1) {
            $p.Range.Delete()
        }
    }
    elseif ($t -eq "Accomplished") {
        # This heading shouldn't be here at all - drop it entirely.
        $p.Range.Delete()
    }
    elseif ($t -eq "Code Executions") {
        $seenCodeExecutions = $seenCodeExecutions + 1
        # Keep the later (2nd-seen-while-scanning-backwards = 1st in doc
        # order) occurrence, drop the earlier duplicate.
        if ($seenCodeExecutions -eq 2) {
            $p.Range.Delete()
        }
    }
    elseif ($t -eq "To-Do List") {
        $seenToDoList = $seenToDoList + 1
        if ($seenToDoList -eq 2) {
            $p.Range.Delete()
        }
    }
}

# --- Update the recorded commit hash for the logging-notebook entry -------
$d.Content.Find.Execute(
    "1c946a2617b405bb34b13a9041bbd6e59a45fdd5",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "0e669cdf6bf8728bf6b0016d38053fc8cb8362e8",
    2) | Out-Null
